$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New build entry + category pages smoke-test rows (endpoint check / page sweep)
$rows = @(
    @("fe11f8d1ec built at 2020-09-10 13:02`n", "CategoryPages"),
    @("test", "ComparePages"),
    @("test", "DealerPages"),
    @("test", "ErrorPages"),
    @("test", "FCVPages"),
    @("test", "HomeOffersPages"),
    @("test", "LCertifiedPages"),
    @("test", "MiscPages"),
    @("test", "ModelPagesAccessories"),
    @("test", "ModelPagesDesign"),
    @("test", "ModelPagesFeatures"),
    @("test", "ModelPagesGallery"),
    @("test", "ModelPagesOffers"),
    @("test", "ModelPagesOverview"),
    @("test", "ModelPagesOwners"),
    @("test", "ModelPagesPackages"),
    @("test", "ModelPagesPerformance"),
    @("test", "ModelPagesSafety"),
    @("test", "ModelPagesSpecifications"),
    @("test", "ModelPagesTechnology"),
    @("test", "SponsoredAthletes")
)

$r = 4
foreach ($row in $rows) {
    # Column B first so new shared-string entries land in the same order
    # Excel originally recorded them (page name before the build/test label).
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $r = $r + 1
}

# Row 4 col A holds a build string with an embedded line break; Excel
# auto-raises that row's height when the value is assigned. Re-autofit it
# back down to the sheet's default (matches the other "built at" row).
$ws.Rows.Item(4).AutoFit()

$ws.Range("A4:B24").Select()
